$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 7892571.428571429
$ws.Range("C11").Value = 9943000
$ws.Range("C13").Value = 10689428.57142857
$ws.Range("C14").Value = -3584000
$ws.Range("C15").Value = 45505333.33333333
